# Fruta / hortaliza, semanal
# Insert a new weekly record at row 213 (pushing the existing rows 213-269
# down to 214-270) and populate it with the new week's data. The new row is
# a copy of the former row 213 except for the date (D) and volume (M)
# columns, which carry the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(213).Insert()

$ws.Cells.Item(213, 1).Value  = 10
$ws.Cells.Item(213, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(213, 3).Value  = "La Araucanía"
$ws.Cells.Item(213, 4).Value  = 44855
$ws.Cells.Item(213, 5).Value  = 9
$ws.Cells.Item(213, 6).Value  = "Fruta"
$ws.Cells.Item(213, 7).Value  = 100101
$ws.Cells.Item(213, 8).Value  = "Berries"
$ws.Cells.Item(213, 9).Value  = 100112025
$ws.Cells.Item(213, 10).Value = "Frutilla"
$ws.Cells.Item(213, 11).Value = "Sin especificar"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 1850
$ws.Cells.Item(213, 14).Value = 10000
$ws.Cells.Item(213, 15).Value = 10000
$ws.Cells.Item(213, 16).Value = 10000
$ws.Cells.Item(213, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(213, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(213, 19).Value = 1429
$ws.Cells.Item(213, 20).Value = 7
